# Update coin Price (column D) and Volume(1h) (column E) figures
# pulled from the latest symbol-list refresh. Each target cell is a
# text-formatted (non-numeric-typed) cell in the source sheet, so we
# force text entry ("@" number format) for the write and then clear
# the format again so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "304.28"
    "E2" = "0.25%"
    "D3" = "37.12"
    "E3" = "3.42%"
    "D4" = "5.041"
    "E4" = "-0.88%"
    "E5" = "-0.20%"
    "D6" = "2.206"
    "E6" = "-3.56%"
    "D7" = "8.014"
    "E7" = "-0.70%"
    "D8" = "0.9235"
    "E8" = "-0.81%"
    "D9" = "0.09870"
    "E9" = "-2.27%"
    "D10" = "0.1895"
    "E10" = "3.63%"
    "D11" = "0.08651"
    "E11" = "0.97%"
    "D12" = "0.03673"
    "E12" = "8.52%"
    "D13" = "0.09941"
    "E13" = "0.41%"
    "D14" = "0.001480"
    "E14" = "-0.87%"
    "D15" = "0.005641"
    "E15" = "-0.43%"
    "D16" = "3.454"
    "E16" = "-0.88%"
    "D17" = "4.025"
    "E17" = "1.21%"
    "D18" = "2.254"
    "E18" = "11.40%"
    "D19" = "0.3412"
    "D20" = "0.1311"
    "E20" = "-1.01%"
    "D21" = "4.763"
    "E21" = "5.06%"
    "E22" = "-0.31%"
    "D23" = "0.04597"
    "E23" = "-0.64%"
    "D24" = "0.001249"
    "E24" = "2.84%"
    "D25" = "0.004483"
    "E25" = "-0.22%"
    "D26" = "0.0001400"
    "E26" = "7.92%"
    "D27" = "0.0002721"
    "E27" = "-19.62%"
    "D39" = "0.01843"
    "E39" = "5.22%"
    "D40" = "0.04786"
    "E40" = "1.23%"
    "D41" = "0.008006"
    "E41" = "2.45%"
    "D42" = "0.1401"
    "E42" = "-1.14%"
    "D43" = "0.007570"
    "D44" = "0.002221"
    "E44" = "0.26%"
    "D45" = "0.01039"
    "E45" = "13.39%"
    "D46" = "0.00006296"
    "E46" = "3.98%"
    "D47" = "0.00000000751"
    "E47" = "0.34%"
    "D48" = "0.0005809"
    "E48" = "0.14%"
    "D49" = "37.79"
    "E49" = "868.54%"
    "D50" = "0.002693"
    "E50" = "0.50%"
    "D51" = "0.00002103"
    "E51" = "0.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
